$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 121733.9523964128
$ws.Range("E5").Value = -0.0322340928438171
$ws.Range("F5").Value = 0.2134466353040298
$ws.Range("G5").Value = -0.5217333196967884
$ws.Range("H5").Value = 6.88967694975133

$ws.Range("D6").Value = 122321.3953810417
$ws.Range("E6").Value = -0.04064292778151905
$ws.Range("F6").Value = 0.2496660240479755
$ws.Range("G6").Value = -1.222091290184444
$ws.Range("H6").Value = 11.17698682676608

$ws.Range("D7").Value = 123695.0447844893
$ws.Range("E7").Value = -0.06155936897188958
$ws.Range("F7").Value = 0.348242544962401
$ws.Range("G7").Value = -1.798899944127958
$ws.Range("H7").Value = 9.709169778814468

$ws.Range("D8").Value = 123998.575986155
$ws.Range("E8").Value = -0.0548853577406161
$ws.Range("F8").Value = 0.217301756108597
$ws.Range("G8").Value = -0.8609034449396482
$ws.Range("H8").Value = 6.879784059782924

$ws.Range("D9").Value = 126239.9451784715
$ws.Range("E9").Value = -0.08156730034373043
$ws.Range("F9").Value = 0.352156707720064
$ws.Range("G9").Value = -1.788058622705788
$ws.Range("H9").Value = 11.57987968636653

$ws.Range("D10").Value = 127371.5404308243
$ws.Range("E10").Value = -0.1152899757433591
$ws.Range("F10").Value = 0.4264961909939701
$ws.Range("G10").Value = -1.901264920955966
$ws.Range("H10").Value = 9.940567024820911

$ws.Range("D11").Value = 129296.5741108983
$ws.Range("E11").Value = -0.1904950759137007
$ws.Range("F11").Value = 0.7499071211480109
$ws.Range("G11").Value = -2.645733907725366
$ws.Range("H11").Value = 13.11949373985481

$ws.Range("D12").Value = 121188.6162944646
$ws.Range("E12").Value = 0.06260678385525156
$ws.Range("F12").Value = 0.04268309918069323
$ws.Range("G12").Value = 2.171951195487561
$ws.Range("H12").Value = 10.60733641771266

$ws.Range("D14").Value = 119374.5170199741
$ws.Range("E14").Value = -0.04994793125251919
$ws.Range("F14").Value = 0.129117363698701
$ws.Range("G14").Value = -1.093875961326501
$ws.Range("H14").Value = 9.425963751313164

$ws.Range("D16").Value = 119426.5474815322
$ws.Range("E16").Value = -0.06858953767293566
$ws.Range("F16").Value = 0.1597017843689158
$ws.Range("G16").Value = -0.810140780766792
$ws.Range("H16").Value = 6.829807266412089

$ws.Range("D17").Value = 119451.2492995463
$ws.Range("E17").Value = -0.05753714113305353
$ws.Range("F17").Value = 0.1332414227260378
$ws.Range("G17").Value = -0.8949883767648086
$ws.Range("H17").Value = 8.679488017370183

